$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: bump volume/number and shift the reporting week by one week ---
# (edit the rightmost run first so the earlier runs character offsets stay valid)
$ws.Range("A8").Characters(21, 2).Text = "11"
$ws.Range("C9").Characters(46, 9).Text = "3/17/2024"
$ws.Range("C9").Characters(27, 8).Text = "3/11/2024"

# --- Crime Complaints table: refreshed weekly figures ---
# Row 16
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -71.428571428571
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = -6.25
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -6.25
$ws.Range("N16").Value = -76.923076923076

# Row 17
$ws.Range("D17").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -63.636363636363
$ws.Range("I17").Value = 16
$ws.Range("K17").Value = -20
$ws.Range("L17").Value = 6.666666666666
$ws.Range("M17").Value = 45.454545454545
$ws.Range("N17").Value = -68

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("J14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("D18").Value = 2
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 266.666666666667
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -5.263157894736
$ws.Range("L18").Value = -30.769230769230
$ws.Range("M18").Value = -28
$ws.Range("N18").Value = -78.048780487804

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = -50
$ws.Range("I19").Value = 33
$ws.Range("J19").Value = 38
$ws.Range("K19").Value = -13.157894736842
$ws.Range("L19").Value = 6.451612903225
$ws.Range("M19").Value = -10.810810810810
$ws.Range("N19").Value = -31.25

# Row 20
$ws.Range("D20").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 12
$ws.Range("K20").Value = 9.090909090909
$ws.Range("L20").Value = 9.090909090909
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -84.615384615384

# Row 21
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 7
$ws.Range("E21").Value = 28.571428571428
$ws.Range("F21").Value = 30
$ws.Range("H21").Value = -23.076923076923
$ws.Range("I21").Value = 94
$ws.Range("J21").Value = 105
$ws.Range("K21").Value = -10.476190476190
$ws.Range("L21").Value = -5.050505050505
$ws.Range("M21").Value = -5.050505050505
$ws.Range("N21").Value = -71.165644171779

# Row 23
$ws.Range("J14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -69.230769230769
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = 63.636363636363

# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = -24.137931034482
$ws.Range("I24").Value = 93
$ws.Range("J24").Value = 127
$ws.Range("K24").Value = -26.771653543307
$ws.Range("L24").Value = 40.909090909090
$ws.Range("M24").Value = 12.048192771084

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -77.777777777777
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -71.428571428571
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 72
$ws.Range("K25").Value = -61.111111111111
$ws.Range("L25").Value = 64.705882352941

# Row 26
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 12
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 24
$ws.Range("J26").Value = 40
$ws.Range("K26").Value = -40
$ws.Range("L26").Value = -31.428571428571
$ws.Range("M26").Value = -48.936170212766

# Row 27
$ws.Range("J14").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 1
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 50

# Row 28
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("J14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("J14").Copy() | Out-Null
$ws.Range("G28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value = 1
$ws.Range("K14").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("H28").Value = -100
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = -66.666666666666
$ws.Range("L28").Value = -66.666666666666

# Row 29
$ws.Range("F29").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("N29").Value = -85.714285714285

# Row 30
$ws.Range("F30").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$ws.Range("N30").Value = -83.333333333333

$excel.CutCopyMode = $false
